$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=N...d=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=3, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=100,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B2").Value = 0.6571428571428571
$ws.Range("C2").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': RobustScaler(), ''model__subsample'': 1.0, ''model__n_estimators'': 100, ''model__max_depth'': 3, ''model__learning_rate'': 0.01, ''model__gamma'': 0, ''model__colsample_bytree'': 0.5}'
$ws.Range("D2").Value = 0.5333333333333333
$ws.Range("E2").Value = '[1 0 0 1 0 0 1 1 0 1 0 0]'
$ws.Range("F2").Value = '[1 1 1 1 1 1 0 1 1 1 1 0]'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.9818333333333333
$ws.Range("I2").Value = 0.005769821471306291
$ws.Range("J2").Value = 0.5867619047619048
$ws.Range("K2").Value = 0.06014813263296432

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()),
                (''selector'',
                 SelectFromModel(estimator=ExtraTreesClassifier(random_state=42))),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.5, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=0, gpu_id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.01,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B3").Value = 0.6095238095238095
$ws.Range("C3").Value = '{''selector'': SelectFromModel(estimator=ExtraTreesClassifier(random_state=42)), ''scaler'': RobustScaler(), ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 7, ''model__learning_rate'': 0.01, ''model__gamma'': 0, ''model__colsample_bytree'': 0.5}'
$ws.Range("D3").Value = 0.5333333333333333
$ws.Range("E3").Value = '[1 0 1 0 0 0 0 1 1 0 1 1]'
$ws.Range("F3").Value = '[1 1 1 1 1 0 1 0 0 1 1 1]'
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.9792142857142858
$ws.Range("I3").Value = 0.005096997445905191
$ws.Range("J3").Value = 0.5362857142857143
$ws.Range("K3").Value = 0.06661386197060812

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', None),
                (''selector'',
                 SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'',
                                                     random_state=42))),
                (''model'',
                 XGBClassifier(base_score=None, booster=None, callbacks=None,
                               colsample_bylevel=None, colsample_bynode=None,
                               colsample_bytree=0.8, early_stopping_rounds=None,
                               enable_categorical=False, eval_metric=None,
                               feature_types=None, gamma=..._id=None,
                               grow_policy=None, importance_type=None,
                               interaction_constraints=None, learning_rate=0.1,
                               max_bin=None, max_cat_threshold=None,
                               max_cat_to_onehot=None, max_delta_step=None,
                               max_depth=7, max_leaves=None,
                               min_child_weight=None, missing=nan,
                               monotone_constraints=None, n_estimators=50,
                               n_jobs=None, num_parallel_tree=None,
                               predictor=None, random_state=42, ...))])'
$ws.Range("B4").Value = 0.6
$ws.Range("C4").Value = '{''selector'': SelectFromModel(estimator=LinearSVC(dual=False, penalty=''l1'', random_state=42)), ''scaler'': None, ''model__subsample'': 0.5, ''model__n_estimators'': 50, ''model__max_depth'': 7, ''model__learning_rate'': 0.1, ''model__gamma'': 0.2, ''model__colsample_bytree'': 0.8}'
$ws.Range("D4").Value = 0.6666666666666666
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[1 0 0 1 1 1 0 0 1 0 1 1]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.9813809523809524
$ws.Range("I4").Value = 0.007273282772671022
$ws.Range("J4").Value = 0.5234285714285715
$ws.Range("K4").Value = 0.09389862041832295
